$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("D2").Value = "'26.898.23"
$ws.Range("E2").Value = "'  -0.68%  "

# Row 3
$ws.Range("D3").Value = "'1.868.83"
$ws.Range("E3").Value = "'  +0.21%  "

# Row 4
$ws.Range("D4").Value = "'1.001"
$ws.Range("E4").Value = "'  -0.04%  "

# Row 5
$ws.Range("D5").Value = "'305.49"
$ws.Range("E5").Value = "'  -0.04%  "

# Row 6
$ws.Range("D6").Value = "'1.000"
$ws.Range("E6").Value = "'  -0.02%  "

# Row 7
$ws.Range("D7").Value = "'0.5080"

# Row 8
$ws.Range("D8").Value = "'0.3667"

# Row 9
$ws.Range("D9").Value = "'0.07181"
$ws.Range("E9").Value = "'  +0.43%  "

# Row 10
$ws.Range("D10").Value = "'0.8893"
$ws.Range("E10").Value = "'  +0.15%  "

# Row 11
$ws.Range("D11").Value = "'20.61"
$ws.Range("E11").Value = "'  -0.34%  "

# Row 12
$ws.Range("D12").Value = "'1.879.28"
$ws.Range("E12").Value = "'  +0.70%  "

# Row 13
$ws.Range("D13").Value = "'0.07490"
$ws.Range("E13").Value = "'  -0.72%  "

# Row 14
$ws.Range("D14").Value = "'94.30"
$ws.Range("E14").Value = "'  +5.38%  "

# Row 15
$ws.Range("D15").Value = "'5.219"
$ws.Range("E15").Value = "'  -1.62%  "

# Row 16
$ws.Range("D16").Value = "'1.000"
$ws.Range("E16").Value = "'  -0.11%  "

# Row 17
$ws.Range("E17").Value = "'  +0.34%  "

# Row 18
$ws.Range("E18").Value = "'  +0.50%  "

# Row 19
$ws.Range("D19").Value = "'0.9997"
$ws.Range("E19").Value = "'  -0.04%  "

# Row 20
$ws.Range("D20").Value = "'26.952.44"
$ws.Range("E20").Value = "'  -0.59%  "

# Row 21
$ws.Range("D21").Value = "'5.012"
$ws.Range("E21").Value = "'  +0.09%  "

# Row 22
$ws.Range("D22").Value = "'2.117.67"
$ws.Range("E22").Value = "'  +1.30%  "

# Row 23
$ws.Range("E23").Value = "'  -1.09%  "

# Row 24
$ws.Range("D24").Value = "'6.382"
$ws.Range("E24").Value = "'  -0.94%  "

# Row 25
$ws.Range("D25").Value = "'147.96"
$ws.Range("E25").Value = "'  +1.96%  "

# Row 26
$ws.Range("D26").Value = "'1.779"
$ws.Range("E26").Value = "'  -3.15%  "

# Row 27
$ws.Range("E27").Value = "'  -0.43%  "

# Row 28
$ws.Range("D28").Value = "'2.076"
$ws.Range("E28").Value = "'  -0.50%  "

# Row 29
$ws.Range("D29").Value = "'113.42"
$ws.Range("E29").Value = "'  +0.55%  "

# Row 30
$ws.Range("D30").Value = "'4.686"
$ws.Range("E30").Value = "'  +0.69%  "

# Row 31
$ws.Range("D31").Value = "'4.712"
$ws.Range("E31").Value = "'  +1.02%  "

# Row 32
$ws.Range("D32").Value = "'0.09140"
$ws.Range("E32").Value = "'  -0.23%  "

# Row 33
$ws.Range("D33").Value = "'0.05044"
$ws.Range("E33").Value = "'  -1.12%  "

# Row 34
$ws.Range("D34").Value = "'0.7508"
$ws.Range("E34").Value = "'  +3.81%  "

# Row 35
$ws.Range("D35").Value = "'2.981"
$ws.Range("E35").Value = "'  -2.90%  "

# Row 36
$ws.Range("E36").Value = "'  -0.16%  "

# Row 37
$ws.Range("D37").Value = "'3.207"
$ws.Range("E37").Value = "'  +3.87%  "

# Row 38
$ws.Range("B38").Value = "TheSandbox"
$ws.Range("C38").Value = "https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand"
$ws.Range("D38").Value = "'0.5633"
$ws.Range("E38").Value = "'  +6.80%  "

# Row 39
$ws.Range("B39").Value = "RenderToken"
$ws.Range("C39").Value = "https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr"
$ws.Range("D39").Value = "'2.498"
$ws.Range("E39").Value = "'  +0.50%  "

# Row 40
$ws.Range("D40").Value = "'0.01990"
$ws.Range("E40").Value = "'  -2.39%  "

# Row 41
$ws.Range("E41").Value = "'  -0.14%  "

# Row 42
$ws.Range("D42").Value = "'6.609"
$ws.Range("E42").Value = "'  +2.06%  "

# Row 43
$ws.Range("D43").Value = "'115.25"
$ws.Range("E43").Value = "'  -0.86%  "

# Row 44
$ws.Range("D44").Value = "'8.546"
$ws.Range("E44").Value = "'  +3.25%  "

# Row 45
$ws.Range("D45").Value = "'0.1486"
$ws.Range("E45").Value = "'  +1.51%  "

# Row 46
$ws.Range("D46").Value = "'0.4770"

# Row 47
$ws.Range("D47").Value = "'0.9998"
$ws.Range("E47").Value = "'  -0.02%  "

# Row 48
$ws.Range("D48").Value = "'10.12"
$ws.Range("E48").Value = "'  +1.59%  "

# Row 49
$ws.Range("E49").Value = "'  -0.39%  "

# Row 50
$ws.Range("D50").Value = "'36.99"
$ws.Range("E50").Value = "'  +1.21%  "

# Row 51
$ws.Range("D51").Value = "'63.11"
$ws.Range("E51").Value = "'  -0.45%  "
